$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("G2").Value = [double]"27.73790633333333"
$ws.Range("H2").Value = [double]"83.213719"
$ws.Range("I2").Value = [double]"0.005442473085408622"
$ws.Range("J2").Value = [double]"0.005456294496964538"
$ws.Range("K2").Value = [double]"3"
$ws.Range("M2").Value = [double]"97.31877133333334"
$ws.Range("N2").Value = [double]"291.956314"
$ws.Range("O2").Value = [double]"0.7059334474701922"
$ws.Range("P2").Value = [double]"0.7149546880603714"
$ws.Range("Q2").Value = [double]"2699.418963719085"
$ws.Range("R2").Value = [double]"24294.77067347177"
$ws.Range("S2").Value = [double]"0.003842023787946242"
$ws.Range("T2").Value = [double]"0.003901003330042802"

$ws.Range("E3").Value = [double]"3"
$ws.Range("G3").Value = [double]"27.73790633333333"
$ws.Range("H3").Value = [double]"83.213719"
$ws.Range("I3").Value = [double]"0.005442473085408622"
$ws.Range("J3").Value = [double]"0.005456294496964538"
$ws.Range("K3").Value = [double]"3"
$ws.Range("M3").Value = [double]"0.9538123333333332"
$ws.Range("N3").Value = [double]"2.861437"
$ws.Range("O3").Value = [double]"0.006918788836773586"
$ws.Range("P3").Value = [double]"0.007007205186661606"
$ws.Range("Q3").Value = [double]"26.4567571615781"
$ws.Range("R3").Value = [double]"238.110814454203"
$ws.Range("S3").Value = [double]"3.765532202776587E-05"
$ws.Range("T3").Value = [double]"3.823337509908309E-05"

$ws.Range("E4").Value = [double]"3"
$ws.Range("G4").Value = [double]"27.73790633333333"
$ws.Range("H4").Value = [double]"83.213719"
$ws.Range("I4").Value = [double]"0.005442473085408622"
$ws.Range("J4").Value = [double]"0.005456294496964538"
$ws.Range("K4").Value = [double]"3"
$ws.Range("M4").Value = [double]"17.45494933333334"
$ws.Range("N4").Value = [double]"52.364848"
$ws.Range("O4").Value = [double]"0.1266151677572303"
$ws.Range("P4").Value = [double]"0.1282332039825957"
$ws.Range("Q4").Value = [double]"484.1637496610791"
$ws.Range("R4").Value = [double]"4357.473746949712"
$ws.Range("S4").Value = [double]"0.0006890996427232234"
$ws.Range("T4").Value = [double]"0.000699678125218368"

$ws.Range("E5").Value = [double]"3"
$ws.Range("G5").Value = [double]"27.73790633333333"
$ws.Range("H5").Value = [double]"83.213719"
$ws.Range("I5").Value = [double]"0.005442473085408622"
$ws.Range("J5").Value = [double]"0.005456294496964538"
$ws.Range("K5").Value = [double]"3"
$ws.Range("M5").Value = [double]"16.91229366666667"
$ws.Range("N5").Value = [double]"50.736881"
$ws.Range("O5").Value = [double]"0.1226788378970111"
$ws.Range("P5").Value = [double]"0.1242465711103312"
$ws.Range("Q5").Value = [double]"469.1116176078265"
$ws.Range("R5").Value = [double]"4222.004558470438"
$ws.Range("S5").Value = [double]"0.0006676762734036904"
$ws.Range("T5").Value = [double]"0.0006779258822160131"

$ws.Range("E6").Value = [double]"3"
$ws.Range("G6").Value = [double]"27.73790633333333"
$ws.Range("H6").Value = [double]"83.213719"
$ws.Range("I6").Value = [double]"0.005442473085408622"
$ws.Range("J6").Value = [double]"0.005456294496964538"
$ws.Range("K6").Value = [double]"2"
$ws.Range("M6").Value = [double]"5.218453999999999"
$ws.Range("N6").Value = [double]"10.436908"
$ws.Range("O6").Value = [double]"0.03785375803879289"
$ws.Range("P6").Value = [double]"0.02555833166004005"
$ws.Range("Q6").Value = [double]"144.7489882568086"
$ws.Range("R6").Value = [double]"868.4939295408519"
$ws.Range("S6").Value = [double]"0.0002060180593077005"
$ws.Range("T6").Value = [double]"0.000139453784388271"

$ws.Range("E7").Value = [double]"3"
$ws.Range("G7").Value = [double]"45.15892033333333"
$ws.Range("H7").Value = [double]"135.476761"
$ws.Range("I7").Value = [double]"0.008860661851212738"
$ws.Range("J7").Value = [double]"0.008883163910879647"
$ws.Range("K7").Value = [double]"3"
$ws.Range("M7").Value = [double]"97.31877133333334"
$ws.Range("N7").Value = [double]"291.956314"
$ws.Range("O7").Value = [double]"0.7059334474701922"
$ws.Range("P7").Value = [double]"0.7149546880603714"
$ws.Range("Q7").Value = [double]"4394.810641579885"
$ws.Range("R7").Value = [double]"39553.29577421896"
$ws.Range("S7").Value = [double]"0.006255037567494223"
$ws.Range("T7").Value = [double]"0.006351059682892107"

$ws.Range("E8").Value = [double]"3"
$ws.Range("G8").Value = [double]"45.15892033333333"
$ws.Range("H8").Value = [double]"135.476761"
$ws.Range("I8").Value = [double]"0.008860661851212738"
$ws.Range("J8").Value = [double]"0.008883163910879647"
$ws.Range("K8").Value = [double]"3"
$ws.Range("M8").Value = [double]"0.9538123333333332"
$ws.Range("N8").Value = [double]"2.861437"
$ws.Range("O8").Value = [double]"0.006918788836773586"
$ws.Range("P8").Value = [double]"0.007007205186661606"
$ws.Range("Q8").Value = [double]"43.07313517395077"
$ws.Range("R8").Value = [double]"387.658216565557"
$ws.Range("S8").Value = [double]"6.130504830259627E-05"
$ws.Range("T8").Value = [double]"6.224615223028106E-05"

$ws.Range("E9").Value = [double]"3"
$ws.Range("G9").Value = [double]"45.15892033333333"
$ws.Range("H9").Value = [double]"135.476761"
$ws.Range("I9").Value = [double]"0.008860661851212738"
$ws.Range("J9").Value = [double]"0.008883163910879647"
$ws.Range("K9").Value = [double]"3"
$ws.Range("M9").Value = [double]"17.45494933333334"
$ws.Range("N9").Value = [double]"52.364848"
$ws.Range("O9").Value = [double]"0.1266151677572303"
$ws.Range("P9").Value = [double]"0.1282332039825957"
$ws.Range("Q9").Value = [double]"788.2466663663699"
$ws.Range("R9").Value = [double]"7094.219997297329"
$ws.Range("S9").Value = [double]"0.001121894186731391"
$ws.Range("T9").Value = [double]"0.001139116569794663"

$ws.Range("E10").Value = [double]"3"
$ws.Range("G10").Value = [double]"45.15892033333333"
$ws.Range("H10").Value = [double]"135.476761"
$ws.Range("I10").Value = [double]"0.008860661851212738"
$ws.Range("J10").Value = [double]"0.008883163910879647"
$ws.Range("K10").Value = [double]"3"
$ws.Range("M10").Value = [double]"16.91229366666667"
$ws.Range("N10").Value = [double]"50.736881"
$ws.Range("O10").Value = [double]"0.1226788378970111"
$ws.Range("P10").Value = [double]"0.1242465711103312"
$ws.Range("Q10").Value = [double]"763.740922346938"
$ws.Range("R10").Value = [double]"6873.668301122441"
$ws.Range("S10").Value = [double]"0.001087015698905158"
$ws.Range("T10").Value = [double]"0.001103702656537836"

$ws.Range("E11").Value = [double]"3"
$ws.Range("G11").Value = [double]"45.15892033333333"
$ws.Range("H11").Value = [double]"135.476761"
$ws.Range("I11").Value = [double]"0.008860661851212738"
$ws.Range("J11").Value = [double]"0.008883163910879647"
$ws.Range("K11").Value = [double]"2"
$ws.Range("M11").Value = [double]"5.218453999999999"
$ws.Range("N11").Value = [double]"10.436908"
$ws.Range("O11").Value = [double]"0.03785375803879289"
$ws.Range("P11").Value = [double]"0.02555833166004005"
$ws.Range("Q11").Value = [double]"235.6597484491647"
$ws.Range("R11").Value = [double]"1413.958490694988"
$ws.Range("S11").Value = [double]"0.0003354093497793696"
$ws.Range("T11").Value = [double]"0.0002270388494247604"

$ws.Range("E12").Value = [double]"3"
$ws.Range("G12").Value = [double]"2612.668416333333"
$ws.Range("H12").Value = [double]"7838.005249"
$ws.Range("I12").Value = [double]"0.5126334109760676"
$ws.Range("J12").Value = [double]"0.5139352671798969"
$ws.Range("K12").Value = [double]"3"
$ws.Range("M12").Value = [double]"97.31877133333334"
$ws.Range("N12").Value = [double]"291.956314"
$ws.Range("O12").Value = [double]"0.7059334474701922"
$ws.Range("P12").Value = [double]"0.7149546880603714"
$ws.Range("Q12").Value = [double]"254261.6801789658"
$ws.Range("R12").Value = [double]"2288355.121610692"
$ws.Range("S12").Value = [double]"0.3618850710987392"
$ws.Range("T12").Value = [double]"0.3674404286298268"

$ws.Range("E13").Value = [double]"3"
$ws.Range("G13").Value = [double]"2612.668416333333"
$ws.Range("H13").Value = [double]"7838.005249"
$ws.Range("I13").Value = [double]"0.5126334109760676"
$ws.Range("J13").Value = [double]"0.5139352671798969"
$ws.Range("K13").Value = [double]"3"
$ws.Range("M13").Value = [double]"0.9538123333333332"
$ws.Range("N13").Value = [double]"2.861437"
$ws.Range("O13").Value = [double]"0.006918788836773586"
$ws.Range("P13").Value = [double]"0.007007205186661606"
$ws.Range("Q13").Value = [double]"2491.995358409201"
$ws.Range("R13").Value = [double]"22427.95822568281"
$ws.Range("S13").Value = [double]"0.003546802321218383"
$ws.Range("T13").Value = [double]"0.003601249869791292"

$ws.Range("E14").Value = [double]"3"
$ws.Range("G14").Value = [double]"2612.668416333333"
$ws.Range("H14").Value = [double]"7838.005249"
$ws.Range("I14").Value = [double]"0.5126334109760676"
$ws.Range("J14").Value = [double]"0.5139352671798969"
$ws.Range("K14").Value = [double]"3"
$ws.Range("M14").Value = [double]"17.45494933333334"
$ws.Range("N14").Value = [double]"52.364848"
$ws.Range("O14").Value = [double]"0.1266151677572303"
$ws.Range("P14").Value = [double]"0.1282332039825957"
$ws.Range("Q14").Value = [double]"45603.99483189858"
$ws.Range("R14").Value = [double]"410435.9534870872"
$ws.Range("S14").Value = [double]"0.06490716532869598"
$ws.Range("T14").Value = [double]"0.06590356595012956"

$ws.Range("E15").Value = [double]"3"
$ws.Range("G15").Value = [double]"2612.668416333333"
$ws.Range("H15").Value = [double]"7838.005249"
$ws.Range("I15").Value = [double]"0.5126334109760676"
$ws.Range("J15").Value = [double]"0.5139352671798969"
$ws.Range("K15").Value = [double]"3"
$ws.Range("M15").Value = [double]"16.91229366666667"
$ws.Range("N15").Value = [double]"50.736881"
$ws.Range("O15").Value = [double]"0.1226788378970111"
$ws.Range("P15").Value = [double]"0.1242465711103312"
$ws.Range("Q15").Value = [double]"44186.21551065426"
$ws.Range("R15").Value = [double]"397675.9395958883"
$ws.Range("S15").Value = [double]"0.0628892711257249"
$ws.Range("T15").Value = [double]"0.06385469471977412"

$ws.Range("E16").Value = [double]"3"
$ws.Range("G16").Value = [double]"2612.668416333333"
$ws.Range("H16").Value = [double]"7838.005249"
$ws.Range("I16").Value = [double]"0.5126334109760676"
$ws.Range("J16").Value = [double]"0.5139352671798969"
$ws.Range("K16").Value = [double]"2"
$ws.Range("M16").Value = [double]"5.218453999999999"
$ws.Range("N16").Value = [double]"10.436908"
$ws.Range("O16").Value = [double]"0.03785375803879289"
$ws.Range("P16").Value = [double]"0.02555833166004005"
$ws.Range("Q16").Value = [double]"13634.08994788835"
$ws.Range("R16").Value = [double]"81804.53968733008"
$ws.Range("S16").Value = [double]"0.01940510110168914"
$ws.Range("T16").Value = [double]"0.0131353280103751"

$ws.Range("E17").Value = [double]"3"
$ws.Range("G17").Value = [double]"2372.267130666667"
$ws.Range("H17").Value = [double]"7116.801392"
$ws.Range("I17").Value = [double]"0.4654641145188886"
$ws.Range("J17").Value = [double]"0.4666461821176285"
$ws.Range("K17").Value = [double]"3"
$ws.Range("M17").Value = [double]"97.31877133333334"
$ws.Range("N17").Value = [double]"291.956314"
$ws.Range("O17").Value = [double]"0.7059334474701922"
$ws.Range("P17").Value = [double]"0.7149546880603714"
$ws.Range("Q17").Value = [double]"230866.1224309322"
$ws.Range("R17").Value = [double]"2077795.101878389"
$ws.Range("S17").Value = [double]"0.3285866870359794"
$ws.Range("T17").Value = [double]"0.3336308755704724"

$ws.Range("E18").Value = [double]"3"
$ws.Range("G18").Value = [double]"2372.267130666667"
$ws.Range("H18").Value = [double]"7116.801392"
$ws.Range("I18").Value = [double]"0.4654641145188886"
$ws.Range("J18").Value = [double]"0.4666461821176285"
$ws.Range("K18").Value = [double]"3"
$ws.Range("M18").Value = [double]"0.9538123333333332"
$ws.Range("N18").Value = [double]"2.861437"
$ws.Range("O18").Value = [double]"0.006918788836773586"
$ws.Range("P18").Value = [double]"0.007007205186661606"
$ws.Range("Q18").Value = [double]"2262.697647191144"
$ws.Range("R18").Value = [double]"20364.2788247203"
$ws.Range("S18").Value = [double]"0.003220447919451988"
$ws.Range("T18").Value = [double]"0.003269885547670483"

$ws.Range("E19").Value = [double]"3"
$ws.Range("G19").Value = [double]"2372.267130666667"
$ws.Range("H19").Value = [double]"7116.801392"
$ws.Range("I19").Value = [double]"0.4654641145188886"
$ws.Range("J19").Value = [double]"0.4666461821176285"
$ws.Range("K19").Value = [double]"3"
$ws.Range("M19").Value = [double]"17.45494933333334"
$ws.Range("N19").Value = [double]"52.364848"
$ws.Range("O19").Value = [double]"0.1266151677572303"
$ws.Range("P19").Value = [double]"0.1282332039825957"
$ws.Range("Q19").Value = [double]"41407.80257091871"
$ws.Range("R19").Value = [double]"372670.2231382685"
$ws.Range("S19").Value = [double]"0.05893481694477973"
$ws.Range("T19").Value = [double]"0.05983953505918937"

$ws.Range("E20").Value = [double]"3"
$ws.Range("G20").Value = [double]"2372.267130666667"
$ws.Range("H20").Value = [double]"7116.801392"
$ws.Range("I20").Value = [double]"0.4654641145188886"
$ws.Range("J20").Value = [double]"0.4666461821176285"
$ws.Range("K20").Value = [double]"3"
$ws.Range("M20").Value = [double]"16.91229366666667"
$ws.Range("N20").Value = [double]"50.736881"
$ws.Range("O20").Value = [double]"0.1226788378970111"
$ws.Range("P20").Value = [double]"0.1242465711103312"
$ws.Range("Q20").Value = [double]"40120.47836961538"
$ws.Range("R20").Value = [double]"361084.3053265383"
$ws.Range("S20").Value = [double]"0.05710259665193857"
$ws.Range("T20").Value = [double]"0.05797918804984249"

$ws.Range("E21").Value = [double]"3"
$ws.Range("G21").Value = [double]"2372.267130666667"
$ws.Range("H21").Value = [double]"7116.801392"
$ws.Range("I21").Value = [double]"0.4654641145188886"
$ws.Range("J21").Value = [double]"0.4666461821176285"
$ws.Range("K21").Value = [double]"2"
$ws.Range("M21").Value = [double]"5.218453999999999"
$ws.Range("N21").Value = [double]"10.436908"
$ws.Range("O21").Value = [double]"0.03785375803879289"
$ws.Range("P21").Value = [double]"0.02555833166004005"
$ws.Range("Q21").Value = [double]"12379.56689709599"
$ws.Range("R21").Value = [double]"74277.40138257593"
$ws.Range("S21").Value = [double]"0.01761956596673899"
$ws.Range("T21").Value = [double]"0.0119266978904538"

$ws.Range("E22").Value = [double]"2"
$ws.Range("G22").Value = [double]"38.730512"
$ws.Range("H22").Value = [double]"77.46102399999999"
$ws.Range("I22").Value = [double]"0.00759933956842245"
$ws.Range("J22").Value = [double]"0.005079092294630384"
$ws.Range("K22").Value = [double]"3"
$ws.Range("M22").Value = [double]"97.31877133333334"
$ws.Range("N22").Value = [double]"291.956314"
$ws.Range("O22").Value = [double]"0.7059334474701922"
$ws.Range("P22").Value = [double]"0.7149546880603714"
$ws.Range("Q22").Value = [double]"3769.205840950923"
$ws.Range("R22").Value = [double]"22615.23504570554"
$ws.Range("S22").Value = [double]"0.005364627980033102"
$ws.Range("T22").Value = [double]"0.003631320847137302"

$ws.Range("E23").Value = [double]"2"
$ws.Range("G23").Value = [double]"38.730512"
$ws.Range("H23").Value = [double]"77.46102399999999"
$ws.Range("I23").Value = [double]"0.00759933956842245"
$ws.Range("J23").Value = [double]"0.005079092294630384"
$ws.Range("K23").Value = [double]"3"
$ws.Range("M23").Value = [double]"0.9538123333333332"
$ws.Range("N23").Value = [double]"2.861437"
$ws.Range("O23").Value = [double]"0.006918788836773586"
$ws.Range("P23").Value = [double]"0.007007205186661606"
$ws.Range("Q23").Value = [double]"36.94164002191466"
$ws.Range("R23").Value = [double]"221.649840131488"
$ws.Range("S23").Value = [double]"5.257822577285305E-05"
$ws.Range("T23").Value = [double]"3.559024187046702E-05"

$ws.Range("E24").Value = [double]"2"
$ws.Range("G24").Value = [double]"38.730512"
$ws.Range("H24").Value = [double]"77.46102399999999"
$ws.Range("I24").Value = [double]"0.00759933956842245"
$ws.Range("J24").Value = [double]"0.005079092294630384"
$ws.Range("K24").Value = [double]"3"
$ws.Range("M24").Value = [double]"17.45494933333334"
$ws.Range("N24").Value = [double]"52.364848"
$ws.Range("O24").Value = [double]"0.1266151677572303"
$ws.Range("P24").Value = [double]"0.1282332039825957"
$ws.Range("Q24").Value = [double]"676.0391246140587"
$ws.Range("R24").Value = [double]"4056.234747684352"
$ws.Range("S24").Value = [double]"0.0009621916542999665"
$ws.Range("T24").Value = [double]"0.0006513082782637681"

$ws.Range("E25").Value = [double]"2"
$ws.Range("G25").Value = [double]"38.730512"
$ws.Range("H25").Value = [double]"77.46102399999999"
$ws.Range("I25").Value = [double]"0.00759933956842245"
$ws.Range("J25").Value = [double]"0.005079092294630384"
$ws.Range("K25").Value = [double]"3"
$ws.Range("M25").Value = [double]"16.91229366666667"
$ws.Range("N25").Value = [double]"50.736881"
$ws.Range("O25").Value = [double]"0.1226788378970111"
$ws.Range("P25").Value = [double]"0.1242465711103312"
$ws.Range("Q25").Value = [double]"655.0217928043573"
$ws.Range("R25").Value = [double]"3930.130756826144"
$ws.Range("S25").Value = [double]"0.0009322781470388403"
$ws.Range("T25").Value = [double]"0.0006310598019607292"

$ws.Range("E26").Value = [double]"2"
$ws.Range("G26").Value = [double]"38.730512"
$ws.Range("H26").Value = [double]"77.46102399999999"
$ws.Range("I26").Value = [double]"0.00759933956842245"
$ws.Range("J26").Value = [double]"0.005079092294630384"
$ws.Range("K26").Value = [double]"2"
$ws.Range("M26").Value = [double]"5.218453999999999"
$ws.Range("N26").Value = [double]"10.436908"
$ws.Range("O26").Value = [double]"0.03785375803879289"
$ws.Range("P26").Value = [double]"0.02555833166004005"
$ws.Range("Q26").Value = [double]"202.113395268448"
$ws.Range("R26").Value = [double]"808.4535810737918"
$ws.Range("S26").Value = [double]"0.0002876635612776882"
$ws.Range("T26").Value = [double]"0.0001298131253981172"
